# Auto-generated edit script: apply 2023-04-08 incremental crime-count updates
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 1755
$ws.Range("J3").Value = 1818
$ws.Range("H4").Value = 1690
$ws.Range("J4").Value = 415
$ws.Range("J5").Value = 127
$ws.Range("J6").Value = 2362
$ws.Range("H7").Value = 26003
$ws.Range("J7").Value = 6477

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 18
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 180
$ws.Range("J8").Value = 388
$ws.Range("J9").Value = 42
$ws.Range("J10").Value = 38
$ws.Range("J12").Value = 17
$ws.Range("J14").Value = 23
$ws.Range("J15").Value = 84
$ws.Range("J19").Value = 225
$ws.Range("J20").Value = 137
$ws.Range("J29").Value = 360
$ws.Range("J33").Value = 271
$ws.Range("J36").Value = 96
$ws.Range("J37").Value = 220
$ws.Range("J41").Value = 41
$ws.Range("J42").Value = 257
$ws.Range("J44").Value = 53
$ws.Range("J50").Value = 36
$ws.Range("J51").Value = 85
$ws.Range("J52").Value = 152
$ws.Range("J55").Value = 73
$ws.Range("H63").Value = 239
$ws.Range("J63").Value = 29
$ws.Range("J64").Value = 44
$ws.Range("J65").Value = 167
$ws.Range("J67").Value = 236
$ws.Range("J69").Value = 17
$ws.Range("J70").Value = 11
$ws.Range("J71").Value = 30
$ws.Range("J75").Value = 26
$ws.Range("J76").Value = 97
$ws.Range("J79").Value = 199
$ws.Range("J84").Value = 67
$ws.Range("J85").Value = 295
$ws.Range("J86").Value = 37
$ws.Range("J88").Value = 72
$ws.Range("J89").Value = 73
$ws.Range("J90").Value = 73
$ws.Range("J91").Value = 76
$ws.Range("J97").Value = 43
$ws.Range("H101").Value = 26003
$ws.Range("J101").Value = 6477

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 74
$ws.Range("J3").Value = 118
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 295

# Sheet 4: Norwood Park
$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 17

# Sheet 5: Little Village
$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 58
$ws.Range("J7").Value = 152

# Sheet 7: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 388

# Sheet 9: Auburn Gresham
$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 57
$ws.Range("J3").Value = 58
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 180

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 25
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 73

# Sheet 12: Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 23

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J5").Value = 8
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 220

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 96
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 236

# Sheet 18: South Deering
$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J2").Value = 22
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 67

# Sheet 19: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 47
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 167

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 271

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 106
$ws.Range("J3").Value = 128
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 360

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J6").Value = 83
$ws.Range("J7").Value = 225

# Sheet 27: Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J2").Value = 18
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 53

# Sheet 28: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J2").Value = 9
$ws.Range("J3").Value = 6

# Sheet 29: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 97

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 20
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 66

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 41

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 56
$ws.Range("J3").Value = 49
$ws.Range("J6").Value = 135
$ws.Range("J7").Value = 257

# Sheet 34: Avondale
$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 38

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 73

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 76

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 54
$ws.Range("J3").Value = 71
$ws.Range("J7").Value = 199

# Sheet 43: Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J3").Value = 10
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 44

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 36
$ws.Range("J4").Value = 16
$ws.Range("J7").Value = 137

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 96

# Sheet 54: Brighton Park
$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 84

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 36

# Sheet 61: Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 42

# Sheet 65: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 43

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 11

# Sheet 68: United Center
$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J2").Value = 15
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 72

# Sheet 70: Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 18

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 37

# Sheet 73: Pullman
$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 26

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J6").Value = 29
$ws.Range("J7").Value = 73

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 85

# Sheet 81: Oakland
$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 30

# Sheet 91: Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 17
